$wb = $excel.ActiveWorkbook
foreach ($name in @("LoginData", "RegisterData", "ForgotPassData")) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Activate()
    $excel.ActiveWindow.Zoom = 110
}
